$wb = $excel.ActiveWorkbook

# --- Sheet "Summary" ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.3772241992882562
$ws1.Range("C2").Value = 0.0625
$ws1.Range("D2").Value = 0.8214285714285714
$ws1.Range("E2").Value = 0.1161616161616162
$ws1.Range("F2").Value = 0.2395833333333333
$ws1.Range("G2").Value = 0.5599250936329588
$ws1.Range("H2").Value = 0.7259898341359015
$ws1.Range("I2").Value = 23
$ws1.Range("J2").Value = 345
$ws1.Range("K2").Value = 189
$ws1.Range("L2").Value = 5

# --- Sheet "Classification Report" ---
$ws2 = $wb.Worksheets.Item("Classification Report")

# row 2 - label "0"
$ws2.Range("B2").Value = 0.9742268041237113
$ws2.Range("C2").Value = 0.3539325842696629
$ws2.Range("D2").Value = 0.5192307692307693

# row 3 - label "1"
$ws2.Range("B3").Value = 0.0625
$ws2.Range("C3").Value = 0.8214285714285714
$ws2.Range("D3").Value = 0.1161616161616162

# row 4 - label "accuracy"
$ws2.Range("B4").Value = 0.3772241992882562
$ws2.Range("C4").Value = 0.3772241992882562
$ws2.Range("D4").Value = 0.3772241992882562
$ws2.Range("E4").Value = 0.3772241992882562

# row 5 - label "macro avg"
$ws2.Range("B5").Value = 0.5183634020618557
$ws2.Range("C5").Value = 0.5876805778491172
$ws2.Range("D5").Value = 0.3176961926961927

# row 6 - label "weighted avg"
$ws2.Range("B6").Value = 0.928802692886231
$ws2.Range("C6").Value = 0.3772241992882562
$ws2.Range("D6").Value = 0.4991490320671815

# --- Sheet "Confusion Matrix" ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

# row 2 - Actual 0
$ws3.Range("B2").Value = 189
$ws3.Range("C2").Value = 345

# row 3 - Actual 1
$ws3.Range("B3").Value = 5
$ws3.Range("C3").Value = 23
